$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("C10").Select()
$win.FreezePanes = $true
